$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "veriset_6"

# Replace spaces with underscores in the header labels (row 1)
$ws.Range("C1").Value = "İş_Deneyimi"
$ws.Range("D1").Value = "Harcama_Miktarı"
$ws.Range("E1").Value = "Çalışma_Saati"
$ws.Range("F1").Value = "Mutluluk_Skoru"

# Move the active selection from L16 to J10
$ws.Range("J10").Select()
